$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 44027
$ws.Range("B24").Value = 4
$ws.Range("C24").Value = "Praca nad poprawnymi aplikacjami konsolowymi."

$ws.Range("C24").Select()
